# SE_Est_sv.xlsx update:
#   - Row 1: turn D1 into a plain numeric header (3) instead of a shared-string
#     header, and extend the header row out to column J with the SE labels
#     (adding a new "SE: $\gamma$" label alongside the existing SPF/SCE ones).
#   - Rows 2-4: new data rows (moment names in A:D, SE values in E:J).
#
# NOTE on write order: new text values are interned into the shared-string
# table in first-seen order, so cells are touched in the same order the
# labels first appear in the final sheet (gamma, then DisgATV, FEATV, Var,
# DisgVar) to keep the saved sharedStrings.xml as close as possible to the
# source layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 (headers) ----
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 'SE: $\hat\lambda_{SPF}$(Q)'
$ws.Range("F1").Value = 'SE: $\hat\lambda_{SPF}$(Q)'
$ws.Range("G1").Value = 'SE: $\gamma$'
$ws.Range("H1").Value = 'SE: $\hat\lambda_{SCE}$(M)'
$ws.Range("I1").Value = 'SE: $\hat\lambda_{SCE}$(M)'
$ws.Range("J1").Value = 'SE: $\gamma$'

# New header cells need the same bold/bordered/centered style ("s=1") that
# A1:E1 already carry - copy it over from the existing styled header cell.
$ws.Range("E1").Copy()
$ws.Range("F1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Text labels in columns A-D (written first, in the order each label
# ----  first appears, so the shared-string table comes out in that order)
$ws.Range("A2").Value = 'DisgATV'
$ws.Range("A3").Value = 'FEATV'
$ws.Range("B2").Value = 'Var'
$ws.Range("B3").Value = 'DisgVar'
$ws.Range("C3").Value = 'DisgATV'
$ws.Range("A4").Value = 'FEATV'
$ws.Range("B4").Value = 'DisgVar'
$ws.Range("C4").Value = 'DisgATV'
$ws.Range("D4").Value = 'Var'

# ---- Row 2 numbers ----
$ws.Range("E2").Value = 0.3
$ws.Range("F2").Value = 0.46
$ws.Range("G2").Value = 2.52
$ws.Range("H2").Value = 0.09
$ws.Range("I2").Value = 0.09
$ws.Range("J2").Value = 0.7

# ---- Row 3 numbers ----
$ws.Range("E3").Value = 0.3
$ws.Range("F3").Value = 0.46
$ws.Range("G3").Value = 2.53
$ws.Range("H3").Value = 0.07
$ws.Range("I3").Value = 0.07
$ws.Range("J3").Value = 0.26

# ---- Row 4 numbers ----
$ws.Range("E4").Value = 0.3
$ws.Range("F4").Value = 0.46
$ws.Range("G4").Value = 1.26
$ws.Range("H4").Value = 0.07
$ws.Range("I4").Value = 0.07
$ws.Range("J4").Value = 0.26

# Match the saved selection state in the target file.
$ws.Range("G7").Select()
